# Auto-generated edit script: apply Chocobo_Profits.xlsx cell-value updates
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4066.6191
$ws.Range("I69").Value = 2799.6667
$ws.Range("J69").Value = 4277.778
$ws.Range("K69").Value = 8399.000100000001
$ws.Range("L69").Value = 12833.334
$ws.Range("M69").Value = -7525.000100000001
$ws.Range("N69").Value = -14581.334
$ws.Range("H72").Value = 4066.6191
$ws.Range("I72").Value = 2799.6667
$ws.Range("J72").Value = 4277.778
$ws.Range("K72").Value = 25197.0003
$ws.Range("L72").Value = 38500.002
$ws.Range("M72").Value = -20829.0003
$ws.Range("N72").Value = -47236.002
$ws.Range("H116").Value = 481829.8
$ws.Range("I116").Value = 910594
$ws.Range("K116").Value = 910594
$ws.Range("M116").Value = -907152
$ws.Range("H129").Value = 817.67
$ws.Range("J129").Value = 869.1429000000001
$ws.Range("L129").Value = 2607.4287
$ws.Range("N129").Value = -12607.4287

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 7409662.5
$ws.Range("I88").Value = 11112861
$ws.Range("J88").Value = 3266.6667
$ws.Range("K88").Value = 11112861
$ws.Range("L88").Value = 3266.6667
$ws.Range("M88").Value = -11112455
$ws.Range("N88").Value = -4078.6667
$ws.Range("H91").Value = 7409662.5
$ws.Range("I91").Value = 11112861
$ws.Range("J91").Value = 3266.6667
$ws.Range("K91").Value = 11112861
$ws.Range("L91").Value = 3266.6667
$ws.Range("M91").Value = -11111457
$ws.Range("N91").Value = -6074.6667
$ws.Range("H102").Value = 2366.8462
$ws.Range("I102").Value = 2177
$ws.Range("K102").Value = 2177
$ws.Range("M102").Value = -555

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2538.4614
$ws.Range("I86").Value = 1920
$ws.Range("J86").Value = 2925
$ws.Range("K86").Value = 1920
$ws.Range("L86").Value = 2925
$ws.Range("M86").Value = -797
$ws.Range("N86").Value = -5171
$ws.Range("H89").Value = 2538.4614
$ws.Range("I89").Value = 1920
$ws.Range("J89").Value = 2925
$ws.Range("K89").Value = 9600
$ws.Range("L89").Value = 14625
$ws.Range("M89").Value = -3984
$ws.Range("N89").Value = -25857
$ws.Range("H99").Value = 5081
$ws.Range("I99").Value = 599.5
$ws.Range("J99").Value = 6201.375
$ws.Range("K99").Value = 599.5
$ws.Range("L99").Value = 6201.375
$ws.Range("M99").Value = 898.5
$ws.Range("N99").Value = -9197.375

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1792.2222
$ws.Range("I31").Value = 1013.6667
$ws.Range("J31").Value = 5685
$ws.Range("K31").Value = 1013.6667
$ws.Range("L31").Value = 5685
$ws.Range("M31").Value = -718.6667
$ws.Range("N31").Value = -6275
$ws.Range("H34").Value = 1792.2222
$ws.Range("I34").Value = 1013.6667
$ws.Range("J34").Value = 5685
$ws.Range("K34").Value = 1013.6667
$ws.Range("L34").Value = 5685
$ws.Range("M34").Value = -811.6667
$ws.Range("N34").Value = -6089

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 9080.125
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 9080.125
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 27240.375
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -28612.375
$ws.Range("H65").Value = 9080.125
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 9080.125
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 81721.125
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -88585.125
$ws.Range("H69").Value = 5058.3076
$ws.Range("I69").Value = 850
$ws.Range("J69").Value = 6928.6665
$ws.Range("K69").Value = 2550
$ws.Range("L69").Value = 20785.9995
$ws.Range("M69").Value = -1739
$ws.Range("N69").Value = -22407.9995
$ws.Range("H70").Value = 4000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 12000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -12630
$ws.Range("H72").Value = 5058.3076
$ws.Range("I72").Value = 850
$ws.Range("J72").Value = 6928.6665
$ws.Range("K72").Value = 7650
$ws.Range("L72").Value = 62357.9985
$ws.Range("M72").Value = -3594
$ws.Range("N72").Value = -70469.9985
$ws.Range("H73").Value = 4000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 12000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -14184
$ws.Range("H75").Value = 1250
$ws.Range("I75").Value = 1000
$ws.Range("J75").Value = 1500
$ws.Range("K75").Value = 3000
$ws.Range("L75").Value = 4500
$ws.Range("M75").Value = -2002
$ws.Range("N75").Value = -6496
$ws.Range("H76").Value = 3000
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H78").Value = 1250
$ws.Range("I78").Value = 1000
$ws.Range("J78").Value = 1500
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 13500
$ws.Range("M78").Value = -4008
$ws.Range("N78").Value = -23484
$ws.Range("H79").Value = 3000
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H80").Value = 17212.375
$ws.Range("I80").Value = 10000
$ws.Range("J80").Value = 19616.5
$ws.Range("K80").Value = 30000
$ws.Range("L80").Value = 58849.5
$ws.Range("M80").Value = -29064
$ws.Range("N80").Value = -60721.5
$ws.Range("H81").Value = 2109.4167
$ws.Range("I81").Value = 1004.3333
$ws.Range("J81").Value = 2477.7778
$ws.Range("K81").Value = 3012.9999
$ws.Range("L81").Value = 7433.3334
$ws.Range("M81").Value = -1889.9999
$ws.Range("N81").Value = -9679.3334
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 17212.375
$ws.Range("I83").Value = 10000
$ws.Range("J83").Value = 19616.5
$ws.Range("K83").Value = 90000
$ws.Range("L83").Value = 176548.5
$ws.Range("M83").Value = -85320
$ws.Range("N83").Value = -185908.5
$ws.Range("H84").Value = 2109.4167
$ws.Range("I84").Value = 1004.3333
$ws.Range("J84").Value = 2477.7778
$ws.Range("K84").Value = 9038.9997
$ws.Range("L84").Value = 22300.0002
$ws.Range("M84").Value = -3422.9997
$ws.Range("N84").Value = -33532.00019999999
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 950
$ws.Range("I86").Value = 900
$ws.Range("K86").Value = 2700
$ws.Range("M86").Value = -1514
$ws.Range("H87").Value = 6937.25
$ws.Range("I87").Value = 675
$ws.Range("J87").Value = 13199.5
$ws.Range("K87").Value = 2025
$ws.Range("L87").Value = 39598.5
$ws.Range("M87").Value = -777
$ws.Range("N87").Value = -42094.5
$ws.Range("H88").Value = 4500
$ws.Range("J88").Value = 4500
$ws.Range("L88").Value = 13500
$ws.Range("N88").Value = -14356
$ws.Range("H89").Value = 950
$ws.Range("I89").Value = 900
$ws.Range("K89").Value = 8100
$ws.Range("M89").Value = -2172
$ws.Range("H90").Value = 6937.25
$ws.Range("I90").Value = 675
$ws.Range("J90").Value = 13199.5
$ws.Range("K90").Value = 6075
$ws.Range("L90").Value = 118795.5
$ws.Range("M90").Value = 165
$ws.Range("N90").Value = -131275.5
$ws.Range("H91").Value = 4500
$ws.Range("J91").Value = 4500
$ws.Range("L91").Value = 13500
$ws.Range("N91").Value = -16464
$ws.Range("H131").Value = 799.2083
$ws.Range("I131").Value = 325
$ws.Range("J131").Value = 819.8261
$ws.Range("K131").Value = 975
$ws.Range("L131").Value = 2459.4783
$ws.Range("M131").Value = 4065
$ws.Range("N131").Value = -12539.4783
$ws.Range("H133").Value = 5154.2856
$ws.Range("I133").Value = 5300
$ws.Range("K133").Value = 15900
$ws.Range("M133").Value = -10840

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 28000
$ws.Range("J62").Value = 28000
$ws.Range("L62").Value = 28000
$ws.Range("N62").Value = -29248
$ws.Range("H65").Value = 28000
$ws.Range("J65").Value = 28000
$ws.Range("L65").Value = 84000
$ws.Range("N65").Value = -90240

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 24600
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 24600
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 24600
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -25848
$ws.Range("H65").Value = 24600
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 24600
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 123000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -129240
$ws.Range("H69").Value = 14633
$ws.Range("J69").Value = 14633
$ws.Range("L69").Value = 14633
$ws.Range("N69").Value = -16131
$ws.Range("H72").Value = 14633
$ws.Range("J72").Value = 14633
$ws.Range("L72").Value = 43899
$ws.Range("N72").Value = -51387
$ws.Range("H86").Value = 31500
$ws.Range("J86").Value = 31500
$ws.Range("L86").Value = 31500
$ws.Range("N86").Value = -33746
$ws.Range("H89").Value = 31500
$ws.Range("J89").Value = 31500
$ws.Range("L89").Value = 157500
$ws.Range("N89").Value = -168732
